$wb = $excel.ActiveWorkbook

# --- Update selection on the previously-active sheet (Manage_Slider) -----
$sliderSheet = $wb.Worksheets.Item("Manage_Slider")
$sliderSheet.Activate() | Out-Null
$sliderSheet.Range("B6").Select() | Out-Null

# --- Add the new "Mobile_Slider" sheet as the last tab -------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sliderSheet)
$newSheet.Name = "Mobile_Slider"

# --- Populate the new sheet's data ----------------------------------------
$newSheet.Range("A1").Value = "Table headers"

$newSheet.Range("A2").Value = "Image"
$newSheet.Range("B2").Value = "Status"
$newSheet.Range("C2").Value = "Action"

$newSheet.Range("A3").Value = "Edit button"
$newSheet.Range("B3").Value = "rgba(0, 123, 255, 1)"

$newSheet.Range("A4").Value = "Delete button"
$newSheet.Range("B4").Value = "rgba(220, 53, 69, 1)"

# --- Column widths (closest achievable values to the authored widths) ----
$newSheet.Columns.Item(1).ColumnWidth = 12
$newSheet.Columns.Item(2).ColumnWidth = 18.333333333333336

# --- Page setup -------------------------------------------------------------
$newSheet.PageSetup.Orientation = 1

# --- Selection on the new sheet, which becomes the active tab ------------
$newSheet.Range("L10").Select() | Out-Null
